# Update the "Skills" line on the resume (re-ordered / revised list of
# technologies) and relocate the stray "_GoBack" bookmark that used to sit
# mid-sentence inside that line.

$d = $word.ActiveDocument

$oldSkills = "Responsive Design, REST API, JSON, HTML5, CSS3, Bootstrap, Media Queries, JavaScript, jQuery, Node, Express, Handlebars, Sequelize, React, MySQL, Visual Studio Code, Postman, Git, Heroku, and Photoshop"
$newSkills = "REST, API, JSON, Responsive Design, Bootstrap, Media Queries, React, HTML, CSS, JavaScript, jQuery, Node, Express, Handlebars, Sequelize, MySQL, Visual Studio Code, Postman, Git, Heroku, and Photoshop"

# Replacing the whole run of text also removes the "_GoBack" bookmark that
# used to live inside it (it was anchored between the old "JSO" / "N, "
# runs), matching the diff.
$d.Content.Find.Execute($oldSkills, $true, $false, $false, $false, $false, $true, 1, $false, $newSkills, 2) | Out-Null

# The bookmark re-appears in the (empty) paragraph right after the
# "Education" heading. Find that heading, then bookmark the paragraph that
# immediately follows it.
$eduHeadingIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Education`r") {
        $eduHeadingIndex = $i
        break
    }
}
$targetPara = $d.Paragraphs.Item($eduHeadingIndex + 1)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $targetPara.Range)
